# Logic tree input file updated.
# The sheet holds a flat decision-tree table (Node1 | Relationship | Node2).
# For each of the three "root problem" blocks (rows whose column A is the
# "engine run too hot", "seeing coolant every time" and "burning coolant
# smell" questions) a new leading row is inserted that links the question
# node to the generic "Possible_Problem" node (reusing the same default
# Possible_Problem text already used elsewhere in the sheet). Inserting
# from the bottom-most block upward keeps the row numbers of the
# not-yet-processed blocks stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default text used for "Possible_Problem" answer nodes (same text already
# present on several rows, e.g. C4 / C7 / C8 before this edit).
$defaultPossibleProblem = $ws.Cells.Item(4, 3).Value2
$possibleProblemLabel = $ws.Cells.Item(4, 2).Value2

# The wrapped-text style (style index 1 in the original file) used by every
# "Possible_Problem"/answer cell in column C - grab it from an existing one
# so the new cells match exactly, then (re)apply WrapText explicitly since
# assigning .Style resets a cell's WrapText flag to the style's default.
$ppStyle = $ws.Cells.Item(4, 3).Style

# --- Block 3: "burning coolant smell" (originally rows 13-14) ---
# Insert a new row ahead of the block (old row 13), push the two existing
# rows down to 14-15.
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).Value2 = $ws.Cells.Item(14, 1).Value2
$ws.Cells.Item(13, 2).Value2 = $possibleProblemLabel
$ws.Cells.Item(13, 3).Value2 = $defaultPossibleProblem
$ws.Cells.Item(13, 3).Style = $ppStyle
$ws.Cells.Item(13, 3).WrapText = $true
$ws.Rows.Item(13).RowHeight = 288

# --- Block 2: "seeing coolant every time" (originally rows 11-12) ---
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value2 = $ws.Cells.Item(12, 1).Value2
$ws.Cells.Item(11, 2).Value2 = $possibleProblemLabel
$ws.Cells.Item(11, 3).Value2 = $defaultPossibleProblem
$ws.Cells.Item(11, 3).Style = $ppStyle
$ws.Cells.Item(11, 3).WrapText = $true
$ws.Rows.Item(11).RowHeight = 288

# --- Block 1: "engine run too hot" (originally rows 7-10) ---
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value2 = $ws.Cells.Item(8, 1).Value2
$ws.Cells.Item(7, 2).Value2 = $possibleProblemLabel
$ws.Cells.Item(7, 3).Value2 = $defaultPossibleProblem
$ws.Cells.Item(7, 3).Style = $ppStyle
$ws.Cells.Item(7, 3).WrapText = $true
$ws.Rows.Item(7).RowHeight = 288

# Update the view to match the edited file: C12 selected (matches the
# author's final selection after entering the new rows).
$ws.Range("C12").Select()
